$d = $word.ActiveDocument

# Smart/curly quotes used in the document text.
$lsq = [char]0x2018   # left single quote '
$rsq = [char]0x2019   # right single quote '

# The paragraph to remove entirely (text + its paragraph mark):
#   "All the 'Successful' projects having returns >100%,"
$target = "All the " + $lsq + "Successful" + $rsq + " projects having returns >100%,"

$findRange = $d.Content
$found = $findRange.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Expand the found range to the whole paragraph (wdParagraph = 4) so the
    # paragraph mark is included and the paragraph disappears completely,
    # causing every following paragraph's content to shift up by one slot.
    $null = $findRange.Expand(4)
    $startPos = $findRange.Start
    $null = $findRange.Delete()

    # The "_GoBack" bookmark used to sit at the very end of the document's
    # last edited paragraph ("... : 61%"). After the deletion above it needs
    # to move to the start of the paragraph that now begins right where the
    # deleted paragraph used to be (it now holds the "All 24 'Journalism' ..."
    # text that shifted up into this slot).
    if ($d.Bookmarks.Exists("_GoBack")) {
        $bm = $d.Bookmarks.Item("_GoBack")
        $null = $bm.Delete()
    }

    $newBookmarkRange = $d.Range($startPos, $startPos)
    $null = $d.Bookmarks.Add("_GoBack", $newBookmarkRange)
}
